# #5: insurance, claim, debt, investment done
#
# 1) "具有相當價值之財產" (assets of considerable value): the
#    property_category value "otherbonds" is renamed to "antique".
# 2) "保險" (insurance): header is translated/expanded with the
#    standard trailing metadata columns (property_category, category,
#    date, legislator_name, legislator_id, source_file, index); the old
#    per-row "amount paid" text column (col E) is replaced by the
#    constant property_category value "insurance".
# 3) "事業投資" (business investment): same standard trailing metadata
#    columns are appended; property_category value is "investment".

# Helper: write $text into $addr as literal text, bypassing Excel's
# auto date/number inference (e.g. "2013-12-26" would otherwise be
# silently converted into a date serial + a new date-formatted style).
# Strategy: give the cell a throwaway value first so it keeps its
# current style, build the literal string in a scratch cell via a
# quoted formula (never auto-converted), copy it, and paste-special
# VALUES ONLY onto the target so the destination's existing style
# (from EntireColumn.Insert) is left untouched.
function Set-LiteralText {
    param($ws, $addr, $text)
    $ws.Range($addr).Value = 0
    $ws.Range("ZZ1").Formula = '="' + $text + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $ws.Range("ZZ1").Clear()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "具有相當價值之財產": otherbonds -> antique
# ---------------------------------------------------------------
$wsAsset = $wb.Worksheets.Item("具有相當價值之財產")
$wsAsset.Range("F2").Value = "antique"
$wsAsset.Range("F3").Value = "antique"

# ---------------------------------------------------------------
# Sheet "保險" (insurance)
# ---------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("保險")

# insert 6 new columns after the existing E column
$wsIns.Range("F1:K5").EntireColumn.Insert()

# header row
$wsIns.Range("B1").Value = "company"
$wsIns.Range("C1").Value = "name"
$wsIns.Range("D1").Value = "owner"
$wsIns.Range("E1").Value = "property_category"
$wsIns.Range("F1").Value = "category"
$wsIns.Range("G1").Value = "date"
$wsIns.Range("H1").Value = "legislator_name"
$wsIns.Range("I1").Value = "legislator_id"
$wsIns.Range("J1").Value = "source_file"
$wsIns.Range("K1").Value = "index"

# data rows: column E used to hold the free-text "amount paid" value,
# it now holds the constant property_category "insurance"; columns
# F:K get the standard trailing metadata.
$insRows = @(2, 3, 4, 5)
$insIndex = @(107, 108, 110, 111)
for ($i = 0; $i -lt $insRows.Length; $i++) {
    $r = $insRows[$i]
    $wsIns.Range("E$r").Value = "insurance"
    $wsIns.Range("F$r").Value = "normal"
    Set-LiteralText $wsIns "G$r" "2013-12-26"
    $wsIns.Range("H$r").Value = "丁守中"
    $wsIns.Range("I$r").Value = 515
    $wsIns.Range("J$r").Value = "tmpc7fb1"
    $wsIns.Range("K$r").Value = $insIndex[$i]
}

# ---------------------------------------------------------------
# Sheet "事業投資" (business investment)
# ---------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("事業投資")

# insert 7 new columns after the existing G column
$wsInv.Range("H1:N2").EntireColumn.Insert()

# header row
$wsInv.Range("B1").Value = "owner"
$wsInv.Range("C1").Value = "company"
$wsInv.Range("D1").Value = "address"
$wsInv.Range("E1").Value = "total"
$wsInv.Range("F1").Value = "register_date"
$wsInv.Range("G1").Value = "register_reason"
$wsInv.Range("H1").Value = "property_category"
$wsInv.Range("I1").Value = "category"
$wsInv.Range("J1").Value = "date"
$wsInv.Range("K1").Value = "legislator_name"
$wsInv.Range("L1").Value = "legislator_id"
$wsInv.Range("M1").Value = "source_file"
$wsInv.Range("N1").Value = "index"

# data row 2
$wsInv.Range("H2").Value = "investment"
$wsInv.Range("I2").Value = "normal"
Set-LiteralText $wsInv "J2" "2013-12-26"
$wsInv.Range("K2").Value = "丁守中"
$wsInv.Range("L2").Value = 515
$wsInv.Range("M2").Value = "tmpc7fb1"
$wsInv.Range("N2").Value = 125
